$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    @(0.674938964098174, 0.4261704314606063, 0.180826786245351, 0.8230139591873777, 1.041926982080775, 1.233865740654225, 1.097992935039133, 1.05454966485305, 1.010707186697167, 1.174665098728651, 0.9915038189280895, 1.032568395975954, 1.073564370337879, 1.29634843634297, 1.222957310708584, 1.077218693888152),
    @(0.6441171104461911, 0.4896700135395199, 0.2134211924131378, 0.8341483890899861, 1.055465555685529, 1.195516700226279, 1.133158779166932, 1.056043184458593, 1.023704670438494, 1.145122786467195, 1.032607704228432, 1.036617911138447, 1.08991916321945, 1.256181930348764, 1.253723894334721, 1.079945567570936),
    @(0.5793104818050514, 0.4867226343139558, 0.2566708241273308, 0.8204014402300149, 1.054261323083509, 1.209187532097304, 1.088920153175076, 1.055354881088704, 1.013937967541522, 1.155656083360111, 0.9920283517487014, 1.034558679712147, 1.095240732127239, 1.266284659087497, 1.200932662416647, 1.077087706667803),
    @(0.5536919398332348, 0.4884483000010756, 0.2710922582746357, 0.8293742171902622, 1.069558625296981, 1.195984760904655, 1.123395728536402, 1.057656432559984, 1.030283208296995, 1.146036960518357, 1.029565199518718, 1.040508158031597, 1.113419186938668, 1.249263133608853, 1.232708510846813, 1.07872705680802),
    @(0.5673294387761182, 0.6076510503033095, 0.06636307852937416, 0.8285866013349112, 1.055806660484498, 1.159902840052803, 1.111685563760645, 1.052832095715404, 1.013676351923526, 1.119272686005423, 0.9902224193208965, 1.032200461878069, 1.098418896908027, 1.207786637198979, 1.257142300420107, 1.075008478551248),
    @(0.5321746027977414, 0.6389400510944171, 0.08671628640832378, 0.8356756939797073, 1.072912175272488, 1.1383286337601, 1.154608467914363, 1.055527090408077, 1.0318548727902, 1.102291510531117, 1.037159841472909, 1.037785052968744, 1.118544264009534, 1.177950555766082, 1.290960085999444, 1.075640470881273),
    @(0.5366745109956695, 0.5065669743504568, 0.08608381246519096, 0.7744302918274106, 1.059760590564652, 1.201099919442132, 1.109326486868178, 1.069523855809862, 1.015161159035475, 1.150433591300263, 0.990336500880593, 1.043403951791956, 1.107401539299718, 1.252835500190896, 1.250998030696123, 1.096237460115484),
    @(0.4502154638456338, 0.556076520392226, 0.1059236198199832, 0.7686914758147056, 1.085685785128206, 1.170075159580303, 1.151356886453636, 1.07816183524919, 1.036671976484882, 1.126542437134207, 1.036041202895931, 1.055375549022575, 1.137885000684985, 1.215122584012405, 1.283668811130095, 1.105249064198194),
    @(0.4702501327887972, 0.5632956988270298, 0.1155076254646437, 0.7819745765695939, 1.068328131448411, 1.177979979492275, 1.105806686968175, 1.067198597827071, 1.017059624007728, 1.133422692002327, 0.990581358706739, 1.041845189435226, 1.12178660647591, 1.228293256946802, 1.242135751984562, 1.093401693712098),
    @(0.3934301481375332, 0.5508912505743382, 0.1258229794731973, 0.7735319205615088, 1.094535969228029, 1.172061730762623, 1.14798815287975, 1.076526192783505, 1.04057283796357, 1.12838796518489, 1.035643677792577, 1.054727696723921, 1.153929128012093, 1.217549395141203, 1.280055421069166, 1.102352208888455),
    @(0.4102269369216154, 0.6341206903601967, 0.02458487439929756, 0.7794388229132639, 1.076070035828222, 1.149115069055726, 1.116683247735937, 1.067980153883513, 1.01954566634613, 1.112586325036329, 0.9897633307988506, 1.042611948531186, 1.137597597004142, 1.18832785078368, 1.266852036480633, 1.093421461787918),
    @(0.3285951053771827, 0.6515978742458469, 0.03330498393649274, 0.7689119920065082, 1.104640730598003, 1.133479191477118, 1.163650389413219, 1.078087320268326, 1.045497389127326, 1.100012818353924, 1.039206406359893, 1.055443225620529, 1.171485320654107, 1.166936724304705, 1.307716257265322, 1.102622378347218),
    @(0.2056467488344851, 0.6150382505849697, 0.03248642585510493, 0.7038970011909848, 1.102457172189288, 1.1568921672678, 1.11573803101557, 1.09126323902638, 1.025394808612565, 1.117991994497274, 0.9899215151459679, 1.060060603778776, 1.184216722729051, 1.200482585034327, 1.266381556183599, 1.126734206627765),
    @(0.2336029685724909, 0.6727178804232677, 0.0305489784264249, 0.7344736847637986, 1.119445577384072, 1.12538773295788, 1.164116949566578, 1.089724424030267, 1.051841723456482, 1.094356349189918, 1.039164607376752, 1.063599317778712, 1.191496590273751, 1.159732103398967, 1.309261531899949, 1.121132922998034),
    @(0.2727914238288439, 0.6403956314770581, 0.0747512823137261, 0.7512854004693271, 1.093796726074921, 1.146557700400793, 1.110682131647957, 1.076657447028959, 1.023374238972219, 1.10800846558575, 0.9901845803956173, 1.048342857589417, 1.166761000655648, 1.191849504539064, 1.252316234787269, 1.11171744693004),
    @(0.3134091833571247, 0.6410177905508562, 0.06230039339440884, 0.7756873358210327, 1.107007508063745, 1.137532613982249, 1.158741798833837, 1.075797852948206, 1.045577035929444, 1.102182699320616, 1.037895381079588, 1.052242220500549, 1.172834550246524, 1.178973051921884, 1.300091878485186, 1.109267182709991),
    @(0.1509077165077166, 0.7097457974267267, 0.0087839249075748, 0.74860038294492, 1.109517515244912, 1.118293858986007, 1.118573423575552, 1.077485008374524, 1.027424867608665, 1.089147757050916, 0.989582636529164, 1.049209035181111, 1.198683548120298, 1.148071519109561, 1.272001225765047, 1.106157764472838),
    @(0.1723639901618951, 0.735607700466184, 0.008957292836096182, 0.7711854479245123, 1.128989879925328, 1.101293498993896, 1.167772174530223, 1.077319093124362, 1.055587566511104, 1.07693141268715, 1.04006896587321, 1.054806646306896, 1.210008287976656, 1.177950555766082, 1.316929734957087, 1.101513263899547)
)

$startRow = 2
$startCol = 4  # column D

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
